# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on sheet "Hoja1" (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $wsHoja1.Range("A1")
$text = $cell.Value()
$text = $text.Replace("✅ 1000 Bs = 3.37 = 12984.17 pesos", "✅ 1000 Bs = 3.36 = 12915.46 pesos")
$text = $text.Replace("✅ 12984.17 pesos = 3.36 = 963.23 Bs", "✅ 12915.46 pesos = 3.34 = 958.67 Bs")
$cell.Value = $text

# --- Update the rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 297.49
$wsTasas.Range("O10").Value = 3842.22
$wsTasas.Range("N12").Value = 3867.9
$wsTasas.Range("O12").Value = 287.1
